# Sprint 4 Presentation edit
#
# The presentation's slide 4 contains a SmartArt timeline diagram whose
# nodes hold text like "Sprint 5 (18.12.2020)" / "Sprint 6 (08.01.2021)".
# The Sprint 5 date moved out one sprint (08.01.2021) and the Sprint 6
# date moved out one sprint as well (15.01.2021).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the graphicFrame shape that hosts the SmartArt (timeline) graphic.
$diagramShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasSmartArt) {
        $diagramShape = $candidate
    }
}

$nodes = $diagramShape.SmartArt.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $tr = $node.TextFrame2.TextRange
    $txt = $tr.Text

    if ($txt.StartsWith("Sprint 5 (")) {
        $tr.Text = "Sprint 5 (08.01.2021)"
    }
    elseif ($txt.StartsWith("Sprint 6 (")) {
        $tr.Text = "Sprint 6 (15.01.2021)"
    }
}
